$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = 4416
$ws.Range("C8").Value = "[-0.9943284774199128, -0.28250962495803833, -4.419668912887573]"
$ws.Range("D8").Value = 4.538939756016481
$ws.Range("E8").Value = 3.758346423445415
$ws.Range("F8").Value = 0.8280229801383968
$ws.Range("G8").Value = 1.033683128043664
$ws.Range("H8").Value = -4.419668912887573
$ws.Range("I8").Value = "[2.0530242919921875, -1.0973854064941406, -5.364692687988281]"

# Row 16
$ws.Range("B16").Value = 6017
$ws.Range("C16").Value = "[-0.8119779080152512, -0.011606216430664062, -5.009792447090149]"
$ws.Range("D16").Value = 5.075181099260022
$ws.Range("E16").Value = 5.00786542014799
$ws.Range("F16").Value = 0.9867363000855581
$ws.Range("G16").Value = 0.8120608520084311
$ws.Range("H16").Value = -5.009792447090149
$ws.Range("I16").Value = "[-2.9125900268554688, 0.1522216796875, -8.711753845214844]"

# Row 19
$ws.Range("B19").Value = 6427
$ws.Range("C19").Value = "[-0.2791014313697815, -1.3958299197256565, -5.224152088165283]"
$ws.Range("D19").Value = 5.414610218111335
$ws.Range("E19").Value = 4.929475358221677
$ws.Range("F19").Value = 0.9104026254250159
$ws.Range("G19").Value = 1.423460141273367
$ws.Range("H19").Value = -5.224152088165283
$ws.Range("I19").Value = "[-0.68609619140625, -1.8134956359863281, -2.36322021484375]"

# Row 20
$ws.Range("B20").Value = 6489
$ws.Range("C20").Value = "[0.08745795488357544, -0.23588845133781433, -6.488113760948181]"
$ws.Range("D20").Value = 6.492989483308295
$ws.Range("E20").Value = 6.441561896461293
$ws.Range("F20").Value = 0.9920795210004253
$ws.Range("G20").Value = 0.2515795209212585
$ws.Range("H20").Value = -6.488113760948181
$ws.Range("I20").Value = "[-0.4549713134765625, -0.05165863037109375, -4.110870361328125]"

# Row 21
$ws.Range("B21").Value = 6562
$ws.Range("C21").Value = "[-2.398833990097046, 1.0760656893253326, -11.748897433280945]"
$ws.Range("D21").Value = 12.03947311046226
$ws.Range("E21").Value = 9.926848122515773
$ws.Range("F21").Value = 0.8245251292508291
$ws.Range("G21").Value = 2.629129490874901
$ws.Range("H21").Value = -11.74889743328094
$ws.Range("I21").Value = "[0.9076499938964844, 0.5763206481933594, -2.2060546875]"
